$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Different Voltage Levels"
$ws.Range("C26").Select()
